$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Juni")

# Enter start/end time for row 8 (2021-06-08): 13:00 - 15:30
$ws.Range("D8").Value = 0.54166666666666663
$ws.Range("E8").Value = 0.64583333333333337

# Set the task/category label in column O to "Coding"
$ws.Range("O8").Value = "Coding"

# Update the selection shown on the Juni sheet (active cell moved to G7)
$ws.Range("G7").Select()

$wb.Save()
